$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting existing rows 51-69 down to 52-70.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new record.
$ws.Cells.Item(51, 1).Value = 1
$ws.Cells.Item(51, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(51, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(51, 4).Value = 45120
$ws.Cells.Item(51, 5).Value = 15
$ws.Cells.Item(51, 6).Value = 100112052
$ws.Cells.Item(51, 7).Value = "Albahaca"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 650
$ws.Cells.Item(51, 11).Value = 900
$ws.Cells.Item(51, 12).Value = 1000
$ws.Cells.Item(51, 13).Value = 938
$ws.Cells.Item(51, 14).Value = "$/paquete"
$ws.Cells.Item(51, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(51, 16).Value = 938
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = "Hortaliza"

# Ensure the date cell keeps the existing date-number-format style (style index 2).
$ws.Cells.Item(51, 4).NumberFormat = $ws.Cells.Item(52, 4).NumberFormat()
